$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated team-specific transition matrix probabilities
# (recomputed from newly added team-specific time data)

# Row 2
$ws.Range("B2").Value = 0.2022900763358779
$ws.Range("C2").Value = 0.5687022900763359
$ws.Range("J2").Value = 0.01526717557251908
$ws.Range("P2").Value = 0.1259541984732824
$ws.Range("S2").Value = 0.08778625954198473

# Row 3
$ws.Range("B3").Value = 0.006756756756756757
$ws.Range("C3").Value = 0.02027027027027027
$ws.Range("J3").Value = 0.02702702702702703
$ws.Range("P3").Value = 0.722972972972973
$ws.Range("S3").Value = 0.222972972972973

# Row 4
$ws.Range("J4").Value = 0.08
$ws.Range("P4").Value = 0.64
$ws.Range("S4").Value = 0.28

# Row 6
$ws.Range("B6").Value = 0.03703703703703703
$ws.Range("D6").Value = 0.02777777777777778
$ws.Range("F6").Value = 0.04629629629629629
$ws.Range("J6").Value = 0.1759259259259259
$ws.Range("O6").Value = 0.03240740740740741
$ws.Range("Q6").Value = 0.1666666666666667
$ws.Range("R6").Value = 0.07407407407407407
$ws.Range("S6").Value = 0.4398148148148148

# Row 7
$ws.Range("B7").Value = 0.1216931216931217
$ws.Range("D7").Value = 0.02645502645502645
$ws.Range("F7").Value = 0.05291005291005291
$ws.Range("J7").Value = 0.1481481481481481
$ws.Range("O7").Value = 0.02116402116402116
$ws.Range("Q7").Value = 0.1851851851851852
$ws.Range("R7").Value = 0.04761904761904762
$ws.Range("S7").Value = 0.3968253968253968

# Row 8
$ws.Range("B8").Value = 0.07992202729044834
$ws.Range("D8").Value = 0.01559454191033138
$ws.Range("F8").Value = 0.05847953216374269
$ws.Range("J8").Value = 0.1169590643274854
$ws.Range("O8").Value = 0.02339181286549707
$ws.Range("Q8").Value = 0.1871345029239766
$ws.Range("R8").Value = 0.09941520467836257
$ws.Range("S8").Value = 0.4191033138401559

# Row 9
$ws.Range("B9").Value = 0.09731543624161074
$ws.Range("D9").Value = 0.01006711409395973
$ws.Range("F9").Value = 0.04697986577181208
$ws.Range("J9").Value = 0.1644295302013423
$ws.Range("O9").Value = 0.01006711409395973
$ws.Range("Q9").Value = 0.1912751677852349
$ws.Range("R9").Value = 0.07718120805369127
$ws.Range("S9").Value = 0.4026845637583892

# Row 10
$ws.Range("B10").Value = 0.1049633848657445
$ws.Range("D10").Value = 0.02441008950366151
$ws.Range("E10").Value = 0.0008136696501220504
$ws.Range("F10").Value = 0.07973962571196094
$ws.Range("J10").Value = 0.08868999186330349
$ws.Range("O10").Value = 0.01627339300244101
$ws.Range("Q10").Value = 0.2107404393816111
$ws.Range("R10").Value = 0.08787632221318145
$ws.Range("S10").Value = 0.386493083807974

# Row 11
$ws.Range("G11").Value = 0.1512027491408935
$ws.Range("J11").Value = 0.1030927835051546
$ws.Range("K11").Value = 0.1993127147766323
$ws.Range("L11").Value = 0.5292096219931272
$ws.Range("S11").Value = 0.01718213058419244

# Row 12
$ws.Range("G12").Value = 0.7444444444444445
$ws.Range("J12").Value = 0.15
$ws.Range("L12").Value = 0.02777777777777778
$ws.Range("S12").Value = 0.07777777777777778

# Row 13
$ws.Range("G13").Value = 0.6578947368421053
$ws.Range("J13").Value = 0.2631578947368421
$ws.Range("S13").Value = 0.07894736842105263

# Row 15
$ws.Range("F15").Value = 0.004464285714285714
$ws.Range("H15").Value = 0.2053571428571428
$ws.Range("I15").Value = 0.1071428571428571
$ws.Range("J15").Value = 0.2991071428571428
$ws.Range("K15").Value = 0.05357142857142857
$ws.Range("M15").Value = 0.01785714285714286
$ws.Range("O15").Value = 0.02678571428571428
$ws.Range("S15").Value = 0.2857142857142857

# Row 16
$ws.Range("F16").Value = 0.005917159763313609
$ws.Range("H16").Value = 0.2130177514792899
$ws.Range("I16").Value = 0.07692307692307693
$ws.Range("J16").Value = 0.3609467455621302
$ws.Range("K16").Value = 0.106508875739645
$ws.Range("M16").Value = 0.005917159763313609
$ws.Range("O16").Value = 0.02958579881656805
$ws.Range("S16").Value = 0.2011834319526627

# Row 17
$ws.Range("F17").Value = 0.006237006237006237
$ws.Range("H17").Value = 0.2162162162162162
$ws.Range("I17").Value = 0.1018711018711019
$ws.Range("J17").Value = 0.395010395010395
$ws.Range("K17").Value = 0.06860706860706861
$ws.Range("M17").Value = 0.02286902286902287
$ws.Range("O17").Value = 0.04781704781704782
$ws.Range("S17").Value = 0.1413721413721414

# Row 18
$ws.Range("H18").Value = 0.1796116504854369
$ws.Range("I18").Value = 0.116504854368932
$ws.Range("J18").Value = 0.3640776699029126
$ws.Range("K18").Value = 0.1310679611650485
$ws.Range("M18").Value = 0.01456310679611651
$ws.Range("N18").Value = 0.004854368932038835
$ws.Range("O18").Value = 0.06310679611650485
$ws.Range("S18").Value = 0.1262135922330097

# Row 19
$ws.Range("F19").Value = 0.01831750339213026
$ws.Range("H19").Value = 0.1913161465400271
$ws.Range("I19").Value = 0.1255088195386703
$ws.Range("J19").Value = 0.3371777476255088
$ws.Range("K19").Value = 0.08955223880597014
$ws.Range("M19").Value = 0.0135685210312076
$ws.Range("N19").Value = 0.00203527815468114
$ws.Range("O19").Value = 0.06716417910447761
$ws.Range("S19").Value = 0.155359565807327

